$wb = $excel.ActiveWorkbook

# Rename sheets (keep Sheet1 as-is)
$wb.Worksheets.Item("Open").Name = "Sheet2"
$wb.Worksheets.Item("Approval").Name = "Sheet3"
$wb.Worksheets.Item("Scheduled").Name = "Sheet4"
$wb.Worksheets.Item("Completed").Name = "Sheet5"

Write-Host "done"
